$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all Column D (Price) cells to Text format first, so that numeric-looking
# strings (e.g. "537.89") are preserved verbatim as text instead of being
# auto-converted to floating point numbers by Excel.
$priceCells = @("D2","D3","D5","D6","D7","D8","D11","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D29","D30","D31","D32","D33","D34","D35","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values for each changed cell.
$ws.Range("D2").Value = "60.265.28"
$ws.Range("E2").Value = "  -4.60%  "
$ws.Range("D3").Value = "2.472.08"
$ws.Range("E3").Value = "  -7.83%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "537.89"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").Value = "149.25"
$ws.Range("E6").Value = "  -5.39%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.573"
$ws.Range("E8").Value = "  -2.43%  "
$ws.Range("E9").Value = "  -5.12%  "
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").Value = "5.35"
$ws.Range("E11").Value = "  +3.88%  "
$ws.Range("E12").Value = "  -3.63%  "
$ws.Range("D13").Value = "2.910.43"
$ws.Range("E13").Value = "  -7.71%  "
$ws.Range("D14").Value = "24.65"
$ws.Range("E14").Value = "  -5.65%  "
$ws.Range("D15").Value = "60.089.13"
$ws.Range("E15").Value = "  -4.63%  "
$ws.Range("D16").Value = "0.0000140"
$ws.Range("E16").Value = "  -4.63%  "
$ws.Range("D17").Value = "2.535.70"
$ws.Range("E17").Value = "  -5.36%  "
$ws.Range("D18").Value = "11.29"
$ws.Range("E18").Value = "  -5.57%  "
$ws.Range("D19").Value = "4.39"
$ws.Range("E19").Value = "  -4.01%  "
$ws.Range("D20").Value = "327.91"
$ws.Range("E20").Value = "  -4.59%  "
$ws.Range("D21").Value = "0.974"
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("D22").Value = "5.81"
$ws.Range("E22").Value = "  -7.87%  "
$ws.Range("D23").Value = "0.477"
$ws.Range("E23").Value = "  -5.43%  "
$ws.Range("D24").Value = "61.51"
$ws.Range("E24").Value = "  -3.04%  "
$ws.Range("D25").Value = "0.162"
$ws.Range("E25").Value = "  -3.27%  "
$ws.Range("D26").Value = "0.986"
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("D27").Value = "7.87"
$ws.Range("E27").Value = "  -3.39%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "1.29"
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.85"
$ws.Range("E30").Value = "  -4.00%  "
$ws.Range("D31").Value = "0.0₃0781"
$ws.Range("E31").Value = "  -8.56%  "
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "159.43"
$ws.Range("E33").Value = "  -4.27%  "
$ws.Range("D34").Value = "4.66"
$ws.Range("E34").Value = "  -3.43%  "
$ws.Range("D35").Value = "18.54"
$ws.Range("E35").Value = "  -5.07%  "
$ws.Range("E36").Value = "  -5.24%  "
$ws.Range("E37").Value = "  -3.38%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "5.85"
$ws.Range("E38").Value = "  -5.39%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "315.80"
$ws.Range("E39").Value = "  -7.32%  "
$ws.Range("D40").Value = "0.869"
$ws.Range("E40").Value = "  -7.10%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "3.80"
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "36.78"
$ws.Range("E42").Value = "  -3.32%  "
$ws.Range("D43").Value = "0.996"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "10.85"
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("D45").Value = "19.95"
$ws.Range("E45").Value = "  -3.74%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.0947"
$ws.Range("E46").Value = "  -2.52%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.583"
$ws.Range("E47").Value = "  -5.54%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "0.0527"
$ws.Range("E48").Value = "  -6.29%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "18.86"
$ws.Range("E49").Value = "  -7.46%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.014.03"
$ws.Range("E50").Value = "  -3.21%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "0.0231"
$ws.Range("E51").Value = "  -3.81%  "
